$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# regcntr_id / device_id pairs for the new rows (row 102 .. row 146)
$data = @(
    @(10002,3000121),
    @(10003,3000122),
    @(10004,3000123),
    @(10005,3000124),
    @(10006,3000125),
    @(10007,3000126),
    @(10008,3000127),
    @(10009,3000128),
    @(10010,3000129),
    @(10002,3000130),
    @(10003,3000131),
    @(10004,3000132),
    @(10005,3000133),
    @(10006,3000134),
    @(10007,3000135),
    @(10008,3000136),
    @(10009,3000137),
    @(10010,3000138),
    @(10002,3000139),
    @(10003,3000140),
    @(10004,3000141),
    @(10005,3000142),
    @(10006,3000143),
    @(10007,3000144),
    @(10008,3000145),
    @(10009,3000146),
    @(10010,3000147),
    @(10002,3000148),
    @(10003,3000149),
    @(10004,3000150),
    @(10005,3000151),
    @(10006,3000152),
    @(10007,3000153),
    @(10008,3000154),
    @(10009,3000155),
    @(10010,3000156),
    @(10002,3000157),
    @(10003,3000158),
    @(10004,3000159),
    @(10005,3000160),
    @(10006,3000161),
    @(10007,3000162),
    @(10008,3000163),
    @(10009,3000164),
    @(10010,3000165)
)

$startRow = 102
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $regcntrId = $data[$i][0]
    $deviceId = $data[$i][1]

    $ws.Cells.Item($r, 1).Value = $regcntrId
    $ws.Cells.Item($r, 2).Value = $deviceId
    $ws.Cells.Item($r, 3).Value = "eng"
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = "superadmin"
    $ws.Cells.Item($r, 6).Value = "now()"
    $ws.Cells.Item($r, 7).Value = "now()"
}

# Update the view to match the saved state: scrolled to row 128, with the
# newly added block selected.
$ws.Application.ActiveWindow.ScrollRow = 128
$ws.Range("A102:B146").Select() | Out-Null

# Page setup / print orientation as captured when the workbook was saved.
$ws.PageSetup.Orientation = 1
